$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "751÷7=107, 2" "841÷4=210, 1"
Replace-Text "987÷8=123, 3" "743÷9=82, 5"
Replace-Text "443÷8=55, 3" "373÷5=74, 3"
Replace-Text "891÷5=178, 1" "777÷6=129, 3"
Replace-Text "997÷4=249, 1" "792÷4=198, 0"
Replace-Text "388÷6=64, 4" "825÷2=412, 1"
Replace-Text "277÷7=39, 4" "545÷9=60, 5"
Replace-Text "942÷8=117, 6" "903÷6=150, 3"
Replace-Text "493÷2=246, 1" "291÷7=41, 4"
Replace-Text "431÷9=47, 8" "134÷5=26, 4"
Replace-Text "352÷9=39, 1" "590÷3=196, 2"
Replace-Text "291÷8=36, 3" "641÷7=91, 4"
Replace-Text "819÷8=102, 3" "495÷9=55, 0"
Replace-Text "465÷8=58, 1" "269÷5=53, 4"
Replace-Text "869÷4=217, 1" "903÷2=451, 1"
Replace-Text "452÷8=56, 4" "183÷4=45, 3"
Replace-Text "288÷9=32, 0" "895÷4=223, 3"
Replace-Text "821÷9=91, 2" "475÷6=79, 1"
Replace-Text "700÷6=116, 4" "554÷7=79, 1"
Replace-Text "874÷8=109, 2" "735÷5=147, 0"
Replace-Text "471÷8=58, 7" "689÷7=98, 3"
Replace-Text "618÷5=123, 3" "101÷3=33, 2"
Replace-Text "214÷5=42, 4" "690÷5=138, 0"
Replace-Text "280÷9=31, 1" "546÷8=68, 2"
Replace-Text "238÷6=39, 4" "243÷8=30, 3"

Write-Output "Done"
